$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Append three new work-log rows (42-44) ------------------------------
$newRows = @(
    @{ Row = 42; Text = "Form filling - implementation, javascript, styles"; Hours = 6;  Serial = 40546 },
    @{ Row = 43; Text = "Writing final report and related LaTeX learning";   Hours = 5;  Serial = 40547 },
    @{ Row = 44; Text = "Writing final report, preparing presentation";     Hours = 10; Serial = 40548 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column B: work item description (goes through shared strings table)
    $ws.Cells.Item($row, 2).Value = $r.Text

    # Column C: hours spent
    $ws.Cells.Item($row, 3).Value2 = $r.Hours

    # Column D: date, stored as the same serial number Excel already uses.
    # Write the raw value first, then clone the existing date cell's format
    # (copy/paste-special keeps the same cellXf/style index instead of
    # fabricating a brand-new numFmt like NumberFormat= would).
    $ws.Cells.Item($row, 4).Value2 = $r.Serial
    $ws.Cells.Item(41, 4).Copy() | Out-Null
    $ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Update the view state to match where the author left the cursor -----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C45").Select()

# --- Recalculate so the SUM(C4:C565) total reflects the new hours --------
$excel.Calculate()
